$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.015627265704854
$ws.Range("D2").Value = 1.021298618812509
$ws.Range("E2").Value = 1.01923626467175
$ws.Range("F2").Value = 1.013957925717873
$ws.Range("I2").Value = 1.0260943693443
$ws.Range("J2").Value = 1.020851701933635
$ws.Range("K2").Value = 1.024136360997896
$ws.Range("L2").Value = 1.022080110856899
$ws.Range("M2").Value = 1.016817517528885
$ws.Range("N2").Value = 1.011059810587822
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.016612735642829
$ws.Range("D3").Value = 1.021999029556605
$ws.Range("E3").Value = 1.020168297976856
$ws.Range("F3").Value = 1.015587432157119
$ws.Range("I3").Value = 1.026247891690876
$ws.Range("J3").Value = 1.021472079261596
$ws.Range("K3").Value = 1.024643485249876
$ws.Range("L3").Value = 1.02281777906419
$ws.Range("M3").Value = 1.018249573845233
$ws.Range("N3").Value = 1.01126579472552
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.017249969601603
$ws.Range("D4").Value = 1.022451383047433
$ws.Range("E4").Value = 1.020771328459327
$ws.Range("F4").Value = 1.016641238609949
$ws.Range("I4").Value = 1.0263450370214
$ws.Range("J4").Value = 1.021872521155145
$ws.Range("K4").Value = 1.024970100685548
$ws.Range("L4").Value = 1.023294424925193
$ws.Range("M4").Value = 1.019175166027082
$ws.Range("N4").Value = 1.011398712945591
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.017517760357028
$ws.Range("D5").Value = 1.0226413464364
$ws.Range("E5").Value = 1.021024829041429
$ws.Range("F5").Value = 1.017084123684324
$ws.Range("I5").Value = 1.026385351399956
$ws.Range("J5").Value = 1.022040631833361
$ws.Range("K5").Value = 1.025107043962203
$ws.Range("L5").Value = 1.023494645665991
$ws.Range("M5").Value = 1.019564040920578
$ws.Range("N5").Value = 1.011454503928084
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.017562717620764
$ws.Range("D6").Value = 1.022673230033764
$ws.Range("E6").Value = 1.021067392125564
$ws.Range("F6").Value = 1.017158478309849
$ws.Range("I6").Value = 1.026392089536035
$ws.Range("J6").Value = 1.022068844581641
$ws.Range("K6").Value = 1.025130015884588
$ws.Range("L6").Value = 1.023528254168211
$ws.Range("M6").Value = 1.019629320584908
$ws.Range("N6").Value = 1.011463866331336
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.017253548236204
$ws.Range("D7").Value = 1.022453922157529
$ws.Range("E7").Value = 1.020774715798945
$ws.Range("F7").Value = 1.016647156984593
$ws.Range("I7").Value = 1.026345577769426
$ws.Range("J7").Value = 1.021874768381533
$ws.Range("K7").Value = 1.024971931966003
$ws.Range("L7").Value = 1.023297100917429
$ws.Range("M7").Value = 1.019180363144165
$ws.Range("N7").Value = 1.011399458772466
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.01596039959586
$ws.Range("D8").Value = 1.021535503579901
$ws.Range("E8").Value = 1.019551261762846
$ws.Range("F8").Value = 1.014508751999717
$ws.Range("I8").Value = 1.026146707025956
$ws.Range("J8").Value = 1.02106156529376
$ws.Range("K8").Value = 1.024308062325896
$ws.Range("L8").Value = 1.02232954922826
$ws.Range("M8").Value = 1.017301707381987
$ws.Range("N8").Value = 1.011129499998952
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.013678369669363
$ws.Range("D9").Value = 1.019910567137366
$ws.Range("E9").Value = 1.017394924940004
$ws.Range("F9").Value = 1.01073575878802
$ws.Range("I9").Value = 1.025779481029547
$ws.Range("J9").Value = 1.019621043975133
$ws.Range("K9").Value = 1.023126539774946
$ws.Range("L9").Value = 1.020619413782585
$ws.Range("M9").Value = 1.013983015980325
$ws.Range("N9").Value = 1.010650980606716
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.012154719004923
$ws.Range("D10").Value = 1.018822877203711
$ws.Range("E10").Value = 1.015957036694228
$ws.Range("F10").Value = 1.008216737351717
$ws.Range("I10").Value = 1.025523389400391
$ws.Range("J10").Value = 1.01865558757503
$ws.Range("K10").Value = 1.022330991546367
$ws.Range("L10").Value = 1.019475808851285
$ws.Range("M10").Value = 1.011764630500579
$ws.Range("N10").Value = 1.010330065898335
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.011494403424756
$ws.Range("D11").Value = 1.018350854160147
$ws.Range("E11").Value = 1.01533433015185
$ws.Range("F11").Value = 1.007125000909104
$ws.Range("I11").Value = 1.025409827224835
$ws.Range("J11").Value = 1.018236314868331
$ws.Range("K11").Value = 1.021984643333739
$ws.Range("L11").Value = 1.018979774386908
$ws.Range("M11").Value = 1.010802556745678
$ws.Range("N11").Value = 1.010190653222627
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.011249046711697
$ws.Range("D12").Value = 1.018175366656922
$ws.Range("E12").Value = 1.015103014778385
$ws.Range("F12").Value = 1.006719325234601
$ws.Range("I12").Value = 1.025367243726669
$ws.Range("J12").Value = 1.018080393700372
$ws.Range("K12").Value = 1.021855712924521
$ws.Range("L12").Value = 1.018795397261037
$ws.Range("M12").Value = 1.010444967386697
$ws.Range("N12").Value = 1.010138800656683
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.011301680476307
$ws.Range("D13").Value = 1.018213016432006
$ws.Range("E13").Value = 1.015152633375961
$ws.Range("F13").Value = 1.00680635124519
$ws.Range("I13").Value = 1.025376396196389
$ws.Range("J13").Value = 1.018113847686963
$ws.Range("K13").Value = 1.021883381666507
$ws.Range("L13").Value = 1.01883495257751
$ws.Range("M13").Value = 1.01052168214826
$ws.Range("N13").Value = 1.010149926310863
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.011474123926159
$ws.Range("D14").Value = 1.018336351513591
$ws.Range("E14").Value = 1.015315209849553
$ws.Range("F14").Value = 1.007091470865405
$ws.Range("I14").Value = 1.025406315448455
$ws.Range("J14").Value = 1.018223430146566
$ws.Range("K14").Value = 1.021973991644935
$ws.Range("L14").Value = 1.018964536328407
$ws.Range("M14").Value = 1.010773003079236
$ws.Range("N14").Value = 1.010186368471772
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.011580360579314
$ws.Range("D15").Value = 1.0184123215078
$ws.Range("E15").Value = 1.015415376646284
$ws.Range("F15").Value = 1.007267121530408
$ws.Range("I15").Value = 1.025424696494196
$ws.Range("J15").Value = 1.018290923029888
$ws.Range("K15").Value = 1.022029782166085
$ws.Range("L15").Value = 1.019044360191672
$ws.Range("M15").Value = 1.010927819196013
$ws.Range("N15").Value = 1.010208812607805
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.012198529969277
$ws.Range("D16").Value = 1.018854181781805
$ws.Range("E16").Value = 1.015998361690724
$ws.Range("F16").Value = 1.00828917086738
$ws.Range("I16").Value = 1.025530869844569
$ws.Range("J16").Value = 1.018683387476347
$ws.Range("K16").Value = 1.022353938106507
$ws.Range("L16").Value = 1.01950871113833
$ws.Range("M16").Value = 1.01182844798468
$ws.Range("N16").Value = 1.010339308665199
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.012586139207214
$ws.Range("D17").Value = 1.019131069015768
$ws.Range("E17").Value = 1.016364027837239
$ws.Range("F17").Value = 1.008930006192784
$ws.Range("I17").Value = 1.025596754076338
$ws.Range("J17").Value = 1.01892924189549
$ws.Range("K17").Value = 1.022556771601112
$ws.Range("L17").Value = 1.019799759208541
$ws.Range("M17").Value = 1.012392982422005
$ws.Range("N17").Value = 1.010421043568985
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.012812170709029
$ws.Range("D18").Value = 1.01929247171734
$ws.Range("E18").Value = 1.016577305938878
$ws.Range("F18").Value = 1.009303700264098
$ws.Range("I18").Value = 1.02563492533993
$ws.Range("J18").Value = 1.019072526600311
$ws.Range("K18").Value = 1.022674900479315
$ws.Range("L18").Value = 1.019969441049323
$ws.Range("M18").Value = 1.012722121741473
$ws.Range("N18").Value = 1.010468674318216
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.012889232419143
$ws.Range("D19").Value = 1.019347488745282
$ws.Range("E19").Value = 1.016650026792533
$ws.Range("F19").Value = 1.009431104539254
$ws.Range("I19").Value = 1.02564789700755
$ws.Range("J19").Value = 1.019121362986989
$ws.Range("K19").Value = 1.022715148770319
$ws.Range("L19").Value = 1.020027284357892
$ws.Range("M19").Value = 1.012834325527107
$ws.Range("N19").Value = 1.010484907743521
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.012544558040819
$ws.Range("D20").Value = 1.019101372075946
$ws.Range("E20").Value = 1.016324796232805
$ws.Range("F20").Value = 1.008861260443976
$ws.Range("I20").Value = 1.025589711996092
$ws.Range("J20").Value = 1.018902876255579
$ws.Range("K20").Value = 1.022535028154805
$ws.Range("L20").Value = 1.019768540947809
$ws.Range("M20").Value = 1.012332428172534
$ws.Range("N20").Value = 1.010412278723078
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.011423346015311
$ws.Range("D21").Value = 1.018300036743387
$ws.Range("E21").Value = 1.015267335534445
$ws.Range("F21").Value = 1.007007514621781
$ws.Range("I21").Value = 1.025397516055596
$ws.Range("J21").Value = 1.018191165946046
$ws.Range("K21").Value = 1.021947317030809
$ws.Range("L21").Value = 1.018926380678594
$ws.Range("M21").Value = 1.010699001800913
$ws.Range("N21").Value = 1.010175639056069
$ws.Range("B22").Value = 1.019999999999999
$ws.Range("C22").Value = 1.010717896698703
$ws.Range("D22").Value = 1.017795296326087
$ws.Range("E22").Value = 1.014602384132147
$ws.Range("F22").Value = 1.005841083421459
$ws.Range("I22").Value = 1.025274352220603
$ws.Range("J22").Value = 1.017742616900673
$ws.Range("K22").Value = 1.021576172300572
$ws.Range("L22").Value = 1.018396141071209
$ws.Range("M22").Value = 1.009670654000475
$ws.Range("N22").Value = 1.010026457930942
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.011091916433892
$ws.Range("D23").Value = 1.018062954934312
$ws.Range("E23").Value = 1.014954895642019
$ws.Range("F23").Value = 1.006459519517656
$ws.Range("I23").Value = 1.025339863823819
$ws.Range("J23").Value = 1.017980502780957
$ws.Range("K23").Value = 1.021773077486013
$ws.Range("L23").Value = 1.01867730151683
$ws.Range("M23").Value = 1.010215930814702
$ws.Range("N23").Value = 1.01010557931093
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.012563346947768
$ws.Range("D24").Value = 1.019114791156473
$ws.Range("E24").Value = 1.016342523332951
$ws.Range("F24").Value = 1.008892323978245
$ws.Range("I24").Value = 1.025592894805798
$ws.Range("J24").Value = 1.018914790117718
$ws.Range("K24").Value = 1.022544853639669
$ws.Range("L24").Value = 1.019782647387809
$ws.Range("M24").Value = 1.012359790474761
$ws.Range("N24").Value = 1.010416239315024
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.014268730024076
$ws.Range("D25").Value = 1.020331428553963
$ws.Range("E25").Value = 1.017952446043288
$ws.Range("F25").Value = 1.011711787042979
$ws.Range("I25").Value = 1.025876405945686
$ws.Range("J25").Value = 1.019994351459898
$ws.Range("K25").Value = 1.023433378101385
$ws.Range("L25").Value = 1.021062142149858
$ws.Range("M25").Value = 1.014841994164459
$ws.Range("N25").Value = 1.010775023977325
